$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.142674565315247
$ws.Range("B1").Value = 3.000122547149658
$ws.Range("C1").Value = 3.634899139404297
$ws.Range("D1").Value = 3.774166107177734
$ws.Range("E1").Value = 1.203986406326294
